# Update column F (dSF) values for several rows as per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -3
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 7
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = -3
